$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# Column D holds price text like '1.00' or '63.901.41' that Excel's smart
# cell-entry parsing would otherwise coerce into a number; forcing the
# cell NumberFormat to Text ("@") before the write keeps it a string,
# matching the original inlineStr/shared-string cell type.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.901.41'
$ws.Range("E2").Value = '  +1.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.306.35'
$ws.Range("E3").Value = '  +5.88%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.62'
$ws.Range("E5").Value = '  +0.70%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.85'
$ws.Range("E6").Value = '  +5.58%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.302.82'
$ws.Range("E8").Value = '  +6.00%  '

# Row 10
$ws.Range("E10").Value = '  +3.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.47'
$ws.Range("E11").Value = '  +5.49%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.472'
$ws.Range("E12").Value = '  +2.62%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("E13").Value = '  +1.20%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.90'
$ws.Range("E14").Value = '  +1.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.849.67'
$ws.Range("E15").Value = '  +5.95%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.121'
$ws.Range("E16").Value = '  +1.14%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.306.49'
$ws.Range("E17").Value = '  +5.99%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.971.69'
$ws.Range("E18").Value = '  +1.54%  '

# Row 19
$ws.Range("E19").Value = '  +2.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.72'
$ws.Range("E20").Value = '  +1.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.31'
$ws.Range("E21").Value = '  +1.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.742'
$ws.Range("E22").Value = '  +6.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.02'
$ws.Range("E23").Value = '  +4.44%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.53'
$ws.Range("E24").Value = '  +4.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.16'
$ws.Range("E25").Value = '  -3.82%  '

# Row 27
$ws.Range("E27").Value = '  +2.44%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  +0.68%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.28'
$ws.Range("E29").Value = '  +4.61%  '

# Row 30
$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.10%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +3.92%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.41'
$ws.Range("E32").Value = '  +4.94%  '

# Row 33
$ws.Range("E33").Value = '  -0.17%  '

# Row 34
$ws.Range("E34").Value = '  +0.58%  '

# Row 35
$ws.Range("E35").Value = '  +2.26%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").Value = '  +2.93%  '

# Row 37
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0765'
$ws.Range("E37").Value = '  +7.82%  '

# Row 38
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.36'
$ws.Range("E38").Value = '  +2.73%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0399'
$ws.Range("E39").Value = '  +2.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '429.05'
$ws.Range("E40").Value = '  +1.85%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.83'
$ws.Range("E41").Value = '  +6.26%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.026.41'
$ws.Range("E42").Value = '  +5.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.45'
$ws.Range("E43").Value = '  +2.03%  '

# Row 44
$ws.Range("E44").Value = '  -6.51%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.270'
$ws.Range("E45").Value = '  +1.85%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.24'
$ws.Range("E46").Value = '  +5.10%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.50'
$ws.Range("E47").Value = '  +2.87%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.999'
$ws.Range("E48").Value = '  +0.04%  '

# Row 49
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.34'
$ws.Range("E49").Value = '  +2.95%  '

# Row 50
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.68'
$ws.Range("E50").Value = '  +13.91%  '

# Row 51
$ws.Range("E51").Value = '  +1.88%  '
